$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actual")

# Insert a new first column for the "ID" field, shifting the existing
# Documento/Nombre/Vinculo/Fecha/Hora.../Rango columns one to the right.
$ws.Columns.Item(1).Insert()

# Header row
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Documento"
$ws.Range("C1").Value = "Nombre"
$ws.Range("D1").Value = "Vinculo"
$ws.Range("E1").Value = "Fecha"
$ws.Range("F1").Value = "Hora Escaneo"
$ws.Range("G1").Value = "Hora Entrada"
$ws.Range("H1").Value = "Hora Salida"
$ws.Range("I1").Value = "Rango"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1053868254"
$ws.Range("C2").Value = "Tatiana Pachon"
$ws.Range("D2").Value = "Administrativa"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2024-02-08"
$ws.Range("F2").Value = "14:09:27"
$ws.Range("G2").Value = "14:09:27"
$ws.Range("H2").Value = "14:09:27"
$ws.Range("I2").Value = "Entrada PM"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "10267084"
$ws.Range("C3").Value = "Ruben Lopez"
$ws.Range("D3").Value = "Administrativa"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2024-02-08"
$ws.Range("F3").Value = "14:09:39"
$ws.Range("G3").Value = "14:09:39"
$ws.Range("H3").Value = "14:09:39"
$ws.Range("I3").Value = "Entrada PM"

# Row 4 (new)
$ws.Range("A4").Value = 3
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "1053868254"
$ws.Range("C4").Value = "Tatiana Pachon"
$ws.Range("D4").Value = "Administrativa"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2024-02-08"
$ws.Range("F4").Value = "14:12:10"
$ws.Range("G4").Value = "14:09:27"
$ws.Range("H4").Value = "14:12:10"
$ws.Range("I4").Value = "Entrada PM"

# Row 5 (new)
$ws.Range("A5").Value = 4
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "10267084"
$ws.Range("C5").Value = "Ruben Lopez"
$ws.Range("D5").Value = "Administrativa"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2024-02-08"
$ws.Range("F5").Value = "14:13:08"
$ws.Range("G5").Value = "14:09:39"
$ws.Range("H5").Value = "14:13:08"
$ws.Range("I5").Value = "Entrada PM"
